# Add, process, save data for age regional model
# Adds new rows to "sheet1" (the derived data map) for:
#   - DNK1 region breakdown (after the existing DNK1 ageband row)
#   - ITA1 region breakdown (after the existing ITA1 ageband row)
#   - NYS1 (USA) ageband + region breakdown rows (appended before the
#     "no care-home-deaths" rows at the bottom of the table)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# --- Insert DNK1 / region row right after the DNK1 / ageband row (row 5) ---
$ws.Rows.Item(6).Insert()
$ws.Range("A6").Value = "DNK1"
$ws.Range("B6").Value = "Denmark"
$ws.Range("C6").Value = "region"
$ws.Range("D6").Value = "data/derived/DNK1/DNK1_regions.RDS"
$ws.Range("E6").Value = "yes"

# --- Insert ITA1 / region row right after the ITA1 / ageband row (row 10) ---
$ws.Rows.Item(11).Insert()
$ws.Range("A11").Value = "ITA1"
$ws.Range("B11").Value = "Italy"
$ws.Range("C11").Value = "region"
$ws.Range("D11").Value = "data/derived/ITA1/ITA1_regions.RDS"
$ws.Range("E11").Value = "yes"

# --- Insert the two new NYS1 (USA) rows after the KEN1 row (row 17),
#     before the block of "_nch" (no care-home-deaths) rows ---
$ws.Rows.Item(18).Insert()
$ws.Range("A18").Value = "NYS1"
$ws.Range("B18").Value = "USA"
$ws.Range("C18").Value = "ageband"
$ws.Range("D18").Value = "data/derived/USA/NYS1_agebands.RDS"
$ws.Range("E18").Value = ""

$ws.Rows.Item(19).Insert()
$ws.Range("A19").Value = "NYS1"
$ws.Range("B19").Value = "USA"
$ws.Range("C19").Value = "region"
$ws.Range("D19").Value = "data/derived/USA/NYS1_regions.RDS"
$ws.Range("E19").Value = ""

# --- Update the selected cell to match the saved state ---
$ws.Range("B16").Select()
